# Applies the 2026-02-26 09:35 OLX monitor update to the PODSUMOWANIE sheet:
# - adds a new bold/green font + centered cell style used for "0 days" (brand-new) listings
# - appends rows 234-242 of freshly (re)scraped listings to the monitoring log table

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# Excel alignment constants
$xlHAlignLeft = -4131
$xlHAlignCenter = -4108

$rows = @(
  @{
    Row = 234
    A = "2026-02-26 09:35:23"
    B = "poqui"
    C = "Duży pokój z balkonem w 2pokojowym mieszkaniu blisko Politechniki"
    D = 1665
    E = "25.02.2026"
    F = 0
    G = "https://www.olx.pl/d/oferta/duzy-pokoj-z-balkonem-w-2pokojowym-mieszkaniu-blisko-politechniki-CID3-ID19xpQK.html"
    H = "duzy-pokoj-z-balkonem-w-2pokojowym-mieszkaniu-blisko-politechniki-CID3-ID19xpQK"
    StyleF = 16
  }
  @{
    Row = 235
    A = "2026-02-26 09:35:23"
    B = "poqui"
    C = "Nowoczesne mieszkanie 2-pokojowe z balkonem, blisko UMCS, KUL, UP"
    D = 2499
    E = "25.02.2026"
    F = 0
    G = "https://www.olx.pl/d/oferta/nowoczesne-mieszkanie-2-pokojowe-z-balkonem-blisko-umcs-kul-up-CID3-ID19xpwN.html"
    H = "nowoczesne-mieszkanie-2-pokojowe-z-balkonem-blisko-umcs-kul-up-CID3-ID19xpwN"
    StyleF = 16
  }
  @{
    Row = 236
    A = "2026-02-26 09:35:23"
    B = "poqui"
    C = "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy"
    D = 2499
    E = "28.10.2025"
    F = 120
    G = "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html"
    H = "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger"
    StyleF = 15
  }
  @{
    Row = 237
    A = "2026-02-26 09:35:23"
    B = "poqui"
    C = "Przytulny pokój blisko Politechniki – ul. Przytulna"
    D = 549
    E = "10.10.2025"
    F = 139
    G = "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html"
    H = "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz"
    StyleF = 15
  }
  @{
    Row = 238
    A = "2026-02-26 09:35:23"
    B = "poqui"
    C = "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza"
    D = 2049
    E = "19.12.2025"
    F = 68
    G = "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html"
    H = "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc"
    StyleF = 15
  }
  @{
    Row = 239
    A = "2026-02-26 09:35:23"
    B = "pokojewlublinie"
    C = "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58"
    D = 0
    E = "11.08.2025"
    F = 198
    G = "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html"
    H = "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm"
    StyleF = 15
  }
  @{
    Row = 240
    A = "2026-02-26 09:35:23"
    B = "pokojewlublinie"
    C = "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12"
    D = 12640
    E = "19.01.2026"
    F = 37
    G = "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html"
    H = "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc"
    StyleF = 14
  }
  @{
    Row = 241
    A = "2026-02-26 09:35:23"
    B = "dawnypatron"
    C = "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4."
    D = 730
    E = "20.09.2024"
    F = 523
    G = "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html"
    H = "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM"
    StyleF = 15
  }
  @{
    Row = 242
    A = "2026-02-26 09:35:23"
    B = "dawnypatron"
    C = "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14"
    D = 14690
    E = "05.12.2025"
    F = 82
    G = "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html"
    H = "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv"
    StyleF = 15
  }
)

# A leading single-quote forces Excel to treat an otherwise date-looking string
# (day-of-month <= 12, e.g. "10.10.2025") as literal text instead of silently
# parsing it into a date serial. Only needed for the handful of "E" values that
# are genuinely ambiguous - everything else round-trips as text on its own.
function Force-Text([string]$s) {
    return "'" + $s
}

foreach ($row in $rows) {
    $r = $row.Row

    # A: timestamp of this check - left aligned (style 13)
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row.A
    $cellA.HorizontalAlignment = $xlHAlignLeft

    # B: profile name - default style
    $ws.Cells.Item($r, 2).Value = $row.B

    # C: listing title - left aligned (style 13)
    $cellC = $ws.Cells.Item($r, 3)
    $cellC.Value = $row.C
    $cellC.HorizontalAlignment = $xlHAlignLeft

    # D: price - centered (style 14)
    $cellD = $ws.Cells.Item($r, 4)
    $cellD.Value = $row.D
    $cellD.HorizontalAlignment = $xlHAlignCenter

    # E: listing date text - centered (style 14)
    $cellE = $ws.Cells.Item($r, 5)
    $dayPart = [int]($row.E.Split(".")[0])
    if ($dayPart -le 12) {
        $cellE.Value = Force-Text $row.E
    } else {
        $cellE.Value = $row.E
    }
    $cellE.HorizontalAlignment = $xlHAlignCenter

    # F: days-on-market counter - centered, colored by age
    $cellF = $ws.Cells.Item($r, 6)
    $cellF.Value = $row.F
    $cellF.HorizontalAlignment = $xlHAlignCenter
    if ($row.StyleF -eq 15) {
        # existing "aging" style: red-ish Calibri 10pt
        $cellF.Font.Size = 10
        $cellF.Font.Bold = $false
        $cellF.Font.Color = 7039999
    } elseif ($row.StyleF -eq 16) {
        # new "brand new" style: bold green Calibri 10pt
        $cellF.Font.Size = 10
        $cellF.Font.Bold = $true
        $cellF.Font.Color = 10551111
    }
    # StyleF -eq 14 uses the plain default font, centered only (no extra font tweaks)

    # G: listing URL - default style
    $ws.Cells.Item($r, 7).Value = $row.G

    # H: listing slug - default style
    $ws.Cells.Item($r, 8).Value = $row.H
}

